$wb = $excel.ActiveWorkbook

# The "home" sheet had its intermediate test steps (rows 2-8) removed,
# leaving only the header row plus the final two steps ("select the
# address" / "add address or pick up point"), which shift up to become
# rows 2-3. Deleting the rows also removes the hyperlink that lived on
# the old row 3 (E3).
$ws = $wb.Worksheets.Item("home")
$ws.Activate()
$ws.Rows("2:8").Delete()
$ws.Hyperlinks.Delete()
$ws.Range("C4").Select()

# Focus moves to the "login" sheet, which becomes the active tab.
$ws1 = $wb.Worksheets.Item("login")
$ws1.Activate()
